$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тесты")

# --- Clear obsolete per-test detail rows (2-4): keep only the numbering (col A)
#     and the Модели: list (col L); drop the stale run parameters/notes. ---
$ws.Range("B2:F2").ClearContents()
$ws.Range("H2:I2").ClearContents()

$ws.Range("B3:F3").ClearContents()
$ws.Range("I3").ClearContents()

$ws.Range("B4:F4").ClearContents()
$ws.Range("H4:I4").ClearContents()

# --- Fix the border "rectangle": it used to close at row 31 (bottom edge)
#     while rows 32+ already carried a separate, offset right-hand border,
#     giving a skewed outline. Drop the stray bottom edge on row 31 so the
#     frame is a proper axis-aligned rectangle that closes one row lower. ---
$ws.Range("A31").Borders.Item(9).LineStyle = -4142
$ws.Range("B31:H31").Borders.Item(9).LineStyle = -4142
$ws.Range("I31").Borders.Item(9).LineStyle = -4142

$ws.Range("B32:H32").Borders.Item(9).LineStyle = -4142
$ws.Range("E26").Borders.Item(9).LineStyle = -4142
$ws.Range("E27").Borders.Item(9).LineStyle = -4142

# Row 33's left corner cell no longer carries any formatting at all now
# that the frame closes at row 32 - clear it outright.
$ws.Range("A33").Clear()

# --- Move the live selection to where work continues. ---
$ws.Range("I32").Select()
